# Updating the Base_Map_Closed to the 2019.2 PreScan Version
#
# Apply the data updates that correspond to re-running the simulation /
# regenerating the goal coordinates for the Base_Map_Closed scenario.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# msg.Pose.Position.X (row 2) and msg.Pose.Position.Y (row 3) got refreshed
# values from the new PreScan run.
$ws.Range("B2").Value = 264
$ws.Range("B3").Value = 156

# The active selection in the saved sheet moved up one row, from B5 to B4.
$ws.Range("B4").Select()
